$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for handoff"
#
# The previous report had two source files:
#   - 02f6eebe-0141-4509-b737-5ea24fd2c875.md   (handed off, with a zh-cn /
#     de-de xlf target each)
#   - 0a6bba8c-aa8a-4898-8d17-3b1e03c81a59.md   (handoff failed)
#
# The new report regenerates the handoff: the failed file's row disappears,
# and the surviving source file gets a new guid/hash and fresh handoff
# timestamps.
# ---------------------------------------------------------------------------

$oldMd   = "02f6eebe-0141-4509-b737-5ea24fd2c875.md"
$newMd   = "1cab4836-11d3-460d-bb6f-b163cee0a381.md"

$oldZhXlf = "02f6eebe-0141-4509-b737-5ea24fd2c875.0065ea26dc5f9f17c56a10e8794c8c09e8c38d53.zh-cn.xlf"
$newZhXlf = "1cab4836-11d3-460d-bb6f-b163cee0a381.a93b885445f8ef946d1d7566735d1e5e95cf65d6.zh-cn.xlf"

$oldDeXlf = "02f6eebe-0141-4509-b737-5ea24fd2c875.0065ea26dc5f9f17c56a10e8794c8c09e8c38d53.de-de.xlf"
$newDeXlf = "1cab4836-11d3-460d-bb6f-b163cee0a381.a93b885445f8ef946d1d7566735d1e5e95cf65d6.de-de.xlf"

$oldZhDate = "2016-01-08 11:27:57"
$newZhDate = "2016-01-08 11:28:42"

$oldDeDate = "2016-01-08 11:28:08"
$newDeDate = "2016-01-08 11:28:50"

$cfgName = ".localization-config"

$mdBase  = "https://github.com/OpenLocalizationTest/oltest/blob/85fe9ad8f0ed817dcd5aca0f9a0b62091aedd5d0/e2e/"
$cfgUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/85fe9ad8f0ed817dcd5aca0f9a0b62091aedd5d0/.localization-config"
$zhXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0157c19d818f4ff43b439009a21727e40f6f3a09/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/"
$deXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1b89321fdc853ad7a212617328768e360428d561/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/"

$newMdUrl   = $mdBase + $newMd
$newZhXlfUrl = $zhXlfBase + $newZhXlf
$newDeXlfUrl = $deXlfBase + $newDeXlf

# ---------------------------------------------------------------------------
# Sheet 1: "Overview" -- A=File Name, B=zh-cn, C=de-de
#   row2 = source file (renamed), row3 = old .localization-config row
#   (the old row3, the failed-handoff source file, is removed entirely)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(3).Delete()
$ws1.Range("A2").Value = $newMd

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMd)
$ws1.Hyperlinks.Add($ws1.Range("A3"), $cfgUrl, [Type]::Missing, [Type]::Missing, $cfgName)

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(3).Delete()
$ws2.Range("A2").Value = $newMd
$ws2.Range("C2").Value = $newZhXlf
$ws2.Range("D2").Value = $newZhDate

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMd)
$ws2.Hyperlinks.Add($ws2.Range("C2"), $newZhXlfUrl, [Type]::Missing, [Type]::Missing, $newZhXlf)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $cfgUrl, [Type]::Missing, [Type]::Missing, $cfgName)

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(3).Delete()
$ws3.Range("A2").Value = $newMd
$ws3.Range("C2").Value = $newDeXlf
$ws3.Range("D2").Value = $newDeDate

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMd)
$ws3.Hyperlinks.Add($ws3.Range("C2"), $newDeXlfUrl, [Type]::Missing, [Type]::Missing, $newDeXlf)
$ws3.Hyperlinks.Add($ws3.Range("A3"), $cfgUrl, [Type]::Missing, [Type]::Missing, $cfgName)
